# Apply weekly update: insert a new data row at row 12 (shifting existing
# rows 12-34 down to 13-35), and populate the new row with the latest
# week's data for "Feria Lagunitas de Puerto Montt - Arándano (blue)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 12; Excel shifts rows 12..34
# down to 13..35 and copies formatting from the row above (row 11).
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(12, 3).Value = "Los Lagos"
$ws.Cells.Item(12, 4).Value = "12/19/2022"
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100101
$ws.Cells.Item(12, 8).Value = "Berries"
$ws.Cells.Item(12, 9).Value = 100101001
$ws.Cells.Item(12, 10).Value = "Arándano (blue)"
$ws.Cells.Item(12, 11).Value = "Sin especificar"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 200
$ws.Cells.Item(12, 14).Value = 3800
$ws.Cells.Item(12, 15).Value = 4000
$ws.Cells.Item(12, 16).Value = 3900
$ws.Cells.Item(12, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(12, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(12, 19).Value = 1950
$ws.Cells.Item(12, 20).Value = 2
